$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Step 1: move the current "New" sheet rows (A2:F7) down onto the end of
# --- the "Previously added" sheet (rows 236:241), preserving formatting via copy/paste.
$src = $ws2.Range("A2:F7")
$dst = $ws1.Range("A236:F241")
$src.Copy($dst)

# --- Step 2: hyperlink each moved row's A cell to the URL held in that cell
# --- (matches the convention used throughout the "Previously added" sheet).
for ($i = 236; $i -le 241; $i++) {
    $cell = $ws1.Cells.Item($i, 1)
    $url = $cell.Value()
    $ws1.Hyperlinks.Add($cell, $url)
}

# --- Step 3: remove the now-migrated rows from the "New" sheet, leaving just
# --- the header and two rows that will hold the newly scraped listings.
$ws2.Rows("4:7").Delete()

# --- Step 4: drop the stale hyperlinks left behind on the "New" sheet; we will
# --- recreate them for the two new listings below.
$ws2.Hyperlinks.Delete()

# --- Step 5: write the two newly scraped forest listings into the "New" sheet.
$ws2.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/raunas-pag/ndcjg.html"
$ws2.Range("B2").Value = "116 000 €"
$ws2.Range("C2").Value = "Cēsis un raj."
$ws2.Range("D2").Value = "21 ha."
$ws2.Range("E2").Value = "42760010118,"
$ws2.Range("F2").Value = 45968.40625

$ws2.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/salacgrivas-l-t/lcnpc.html"
$ws2.Range("B3").Value = "90 000 €"
$ws2.Range("C3").Value = "Limbaži un raj."
$ws2.Range("D3").Value = "14 ha."
$ws2.Range("E3").Value = "66720030013"
$ws2.Range("F3").Value = 45967.89583333333

# --- Step 6: hyperlink the two new listings' A cells to their URLs.
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/raunas-pag/ndcjg.html")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/salacgrivas-l-t/lcnpc.html")
